$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "954÷2=" "591÷9="
Replace-Text "976÷4=" "696÷4="
Replace-Text "819÷7=" "433÷2="
Replace-Text "552÷9=" "486÷9="
Replace-Text "826÷4=" "573÷7="
Replace-Text "481÷7=" "186÷4="
Replace-Text "469÷4=" "622÷7="
Replace-Text "374÷7=" "620÷3="
Replace-Text "590÷9=" "596÷3="
Replace-Text "429÷3=" "756÷9="
Replace-Text "816÷5=" "637÷7="
Replace-Text "241÷6=" "837÷2="
Replace-Text "210÷6=" "668÷4="
Replace-Text "517÷9=" "342÷8="
Replace-Text "573÷3=" "666÷2="
Replace-Text "966÷6=" "523÷2="
Replace-Text "435÷3=" "296÷8="
Replace-Text "290÷5=" "576÷2="
Replace-Text "651÷3=" "824÷3="
Replace-Text "422÷5=" "512÷9="
Replace-Text "433÷3=" "474÷6="
Replace-Text "984÷9=" "420÷8="
Replace-Text "545÷3=" "749÷2="
Replace-Text "969÷8=" "382÷4="
Replace-Text "979÷7=" "909÷5="
